# Update "想去人数" (attendance / interest count) figures across the
# workbook's sheets to reflect the latest scrape (gh-pages output
# generated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsShow    = $wb.Worksheets.Item("演出")
$wsAll     = $wb.Worksheets.Item("全部类型")

# --- 展览 (Exhibitions) sheet: column F updates ---
$exhibitUpdates = @{
    "F2"  = 415
    "F3"  = 1029
    "F4"  = 5792
    "F5"  = 540
    "F6"  = 1004
    "F7"  = 1019
    "F8"  = 839
    "F11" = 609
    "F12" = 41
    "F15" = 1961
    "F16" = 1505
    "F17" = 1048
    "F21" = 614
    "F25" = 527
    "F26" = 3379
    "F27" = 190
    "F28" = 126
    "F30" = 150
    "F31" = 47
    "F32" = 469
    "F37" = 315
    "F38" = 802
    "F40" = 66
    "F42" = 80
}
foreach ($cell in $exhibitUpdates.Keys) {
    $wsExhibit.Range($cell).Value = $exhibitUpdates[$cell]
}

# --- 演出 (Shows) sheet: column F updates ---
$showUpdates = @{
    "F4" = 507
    "F6" = 285
}
foreach ($cell in $showUpdates.Keys) {
    $wsShow.Range($cell).Value = $showUpdates[$cell]
}

# --- 全部类型 (All types) sheet: column F updates ---
$allUpdates = @{
    "F3"  = 415
    "F4"  = 1029
    "F6"  = 5792
    "F7"  = 540
    "F8"  = 1004
    "F10" = 507
    "F11" = 1019
    "F12" = 839
    "F14" = 285
    "F17" = 609
    "F18" = 41
    "F22" = 1961
    "F23" = 1505
    "F24" = 1048
    "F29" = 614
    "F32" = 3379
    "F33" = 190
    "F34" = 126
    "F36" = 150
    "F37" = 47
    "F38" = 469
    "F42" = 315
    "F43" = 802
    "F46" = 80
}
foreach ($cell in $allUpdates.Keys) {
    $wsAll.Range($cell).Value = $allUpdates[$cell]
}
